$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 4): rename the repeated C/B/A columns to Z2/Z1/Z0 ---
$ws.Range("K4").Value = "Z2"
$ws.Range("L4").Value = "Z1"
$ws.Range("M4").Value = "Z0"

# Remove the stray "Assuming done is 1 rather than x" note cell
$ws.Range("T4").ClearContents()

# --- Data correction: M7 changes from 1 to 0 ---
$ws.Range("M7").Value = 0

# --- Updated boolean equations in column T ---
$ws.Range("T13").Value = "C+ = CB" + [char]0x2019 + " + CA" + [char]0x2019 + " + C" + [char]0x2019 + "BA"
$ws.Range("T14").Value = "B+ = B" + [char]0x2019 + "A + BA" + [char]0x2019
$ws.Range("T15").Value = "A+ = BA" + [char]0x2019 + " + CA" + [char]0x2019 + "done +A" + [char]0x2019 + "donesensor"

$ws.Range("T17").Value = "Z2 = C"
$ws.Range("T18").Value = "Z1 = B               "
$ws.Range("T19").Value = "Z0 = A               "

$ws.Range("T21").Value = "T8 = CBA"
$ws.Range("T22").Value = "T4 = C" + [char]0x2019 + "BA"

# --- Update the selected range shown when the sheet was last saved ---
$ws.Range("R24").Select()
